$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B21: currently stored as inline string "1" -> should become numeric 1
$ws.Range("B21").Value = 1

# Add new row 22 data
$ws.Range("A22").Value = "Ying Tang"
$ws.Range("B22").Value = "'3"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "misspelled multiple times"
$ws.Range("D22").Value = "DFT"
$ws.Range("E22").Value = "WRI"
$ws.Range("F22").Value = "9386b51e-53f7-4a13-b66e-3217e88401e6"
$ws.Range("G22").Value = "HyRnez-RW_annotated.xlsx"
$ws.Range("H22").Value = '"Krasner" misspelled multiple times as "Kramer"'
